$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.307.40"

$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "2.634.49"

$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.87"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.10"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("B9").Value = "LidoStakedEther"

$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"

$ws.Range("D9").Value = "2.632.89"

$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("B10").Value = "Dogecoin"

$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D14").Value = "3.115.55"

$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "72.182.54"

$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.71"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("D18").Value = "2.635.87"

$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.06"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "375.81"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.85"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.47"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.23"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.40"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("D28").Value = "2.770.01"

$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").Value = "0.0₃0946"

$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "490.24"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -3.77%  "

$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.48"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +8.48%  "

$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.88"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("E40").Value = "  -1.26%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.72"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  -4.74%  "

$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("E44").Value = "  -2.52%  "

$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.00"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.17"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("E48").Value = "  -2.64%  "

$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("E50").Value = "  -3.00%  "

$ws.Range("E51").Value = "  +0.87%  "
